$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$g = $s.Shapes.Item(1)

$it = $g.GroupItems.Item(3)
$it.Left = 272.670000
$it.Top = 126.362441
$it.Width = 163.600472
$it.Height = 141.622756
$it = $g.GroupItems.Item(4)
$it.Left = 293.641811
$it.Top = 126.362441
$it.Width = 0.000000
$it.Height = 141.622756
$it = $g.GroupItems.Item(5)
$it.Left = 335.585512
$it.Top = 126.362441
$it.Width = 0.000000
$it.Height = 141.622756
$it = $g.GroupItems.Item(6)
$it.Left = 377.529213
$it.Top = 126.362441
$it.Width = 0.000000
$it.Height = 141.622756
$it = $g.GroupItems.Item(7)
$it.Left = 419.472913
$it.Top = 126.362441
$it.Width = 0.000000
$it.Height = 141.622756
$it = $g.GroupItems.Item(8)
$it.Left = 272.670000
$it.Top = 247.753386
$it.Width = 163.600472
$it.Height = 0.000000
$it = $g.GroupItems.Item(9)
$it.Left = 272.670000
$it.Top = 214.033701
$it.Width = 163.600472
$it.Height = 0.000000
$it = $g.GroupItems.Item(10)
$it.Left = 272.670000
$it.Top = 180.313937
$it.Width = 163.600472
$it.Height = 0.000000
$it = $g.GroupItems.Item(11)
$it.Left = 272.670000
$it.Top = 146.594252
$it.Width = 163.600472
$it.Height = 0.000000
$it = $g.GroupItems.Item(12)
$it.Left = 272.670000
$it.Top = 126.362441
$it.Width = 0.000000
$it.Height = 141.622756
$it = $g.GroupItems.Item(13)
$it.Left = 314.613701
$it.Top = 126.362441
$it.Width = 0.000000
$it.Height = 141.622756
$it = $g.GroupItems.Item(14)
$it.Left = 356.557402
$it.Top = 126.362441
$it.Width = 0.000000
$it.Height = 141.622756
$it = $g.GroupItems.Item(15)
$it.Left = 398.501024
$it.Top = 126.362441
$it.Width = 0.000000
$it.Height = 141.622756
$it = $g.GroupItems.Item(16)
$it.Left = 272.670000
$it.Top = 131.420394
$it.Width = 17.951890
$it.Height = 30.347717
$it = $g.GroupItems.Item(17)
$it.Left = 272.670000
$it.Top = 165.140079
$it.Width = 122.307717
$it.Height = 30.347717
$it = $g.GroupItems.Item(18)
$it.Left = 272.670000
$it.Top = 198.859764
$it.Width = 45.970236
$it.Height = 30.347717
$it = $g.GroupItems.Item(19)
$it.Left = 272.670000
$it.Top = 232.579528
$it.Width = 160.392598
$it.Height = 30.347717
$it = $g.GroupItems.Item(20)
$it.Left = 377.541890
$it.Top = 176.203150
$it.Width = 4.005276
$it.Height = 5.872520
$it = $g.GroupItems.Item(21)
$it.Left = 382.218268
$it.Top = 176.115669
$it.Width = 4.109528
$it.Height = 5.960000
$it = $g.GroupItems.Item(22)
$it.Left = 386.965433
$it.Top = 176.115669
$it.Width = 4.134488
$it.Height = 6.043386
$it = $g.GroupItems.Item(23)
$it.Left = 388.153307
$it.Top = 177.032598
$it.Width = 1.708819
$it.Height = 2.133937
$it = $g.GroupItems.Item(24)
$it.Left = 415.555906
$it.Top = 243.555118
$it.Width = 4.134488
$it.Height = 6.043386
$it = $g.GroupItems.Item(25)
$it.Left = 416.743780
$it.Top = 244.472047
$it.Width = 1.708819
$it.Height = 2.133937
$it = $g.GroupItems.Item(26)
$it.Left = 420.269764
$it.Top = 243.642598
$it.Width = 4.247008
$it.Height = 5.955906
$it = $g.GroupItems.Item(27)
$it.Left = 425.067008
$it.Top = 243.555118
$it.Width = 4.126142
$it.Height = 6.043386
$it = $g.GroupItems.Item(28)
$it.Left = 426.296535
$it.Top = 246.597638
$it.Width = 1.721260
$it.Height = 2.088031
$it = $g.GroupItems.Item(29)
$it.Left = 294.833071
$it.Top = 143.657953
$it.Width = 3.680236
$it.Height = 5.872520
$it = $g.GroupItems.Item(30)
$it.Left = 299.263543
$it.Top = 143.570472
$it.Width = 4.080315
$it.Height = 6.043386
$it = $g.GroupItems.Item(31)
$it.Left = 300.022126
$it.Top = 144.178976
$it.Width = 2.559055
$it.Height = 4.822205
$it = $g.GroupItems.Item(32)
$it.Left = 304.114961
$it.Top = 143.657953
$it.Width = 3.880236
$it.Height = 5.872520
$it = $g.GroupItems.Item(33)
$it.Left = 322.630551
$it.Top = 211.009843
$it.Width = 3.888583
$it.Height = 5.960000
$it = $g.GroupItems.Item(34)
$it.Left = 327.386142
$it.Top = 211.097402
$it.Width = 3.880236
$it.Height = 5.872520
$it = $g.GroupItems.Item(35)
$it.Left = 331.891654
$it.Top = 211.097402
$it.Width = 4.301181
$it.Height = 5.872520
$it = $g.GroupItems.Item(36)
$it.Left = 332.583465
$it.Top = 211.993071
$it.Width = 2.075591
$it.Height = 3.055433
$it = $g.GroupItems.Item(37)
$it.Left = 272.670000
$it.Top = 126.362441
$it.Width = 163.600472
$it.Height = 141.622756
$it = $g.GroupItems.Item(38)
$it.Left = 180.920472
$it.Top = 244.498504
$it.Width = 4.627717
$it.Height = 6.281969
$it = $g.GroupItems.Item(39)
$it.Left = 181.719764
$it.Top = 245.181732
$it.Width = 2.981969
$it.Height = 4.911260
$it = $g.GroupItems.Item(40)
$it.Left = 186.717008
$it.Top = 249.706299
$it.Width = 1.091339
$it.Height = 1.198819
$it = $g.GroupItems.Item(41)
$it.Left = 191.383386
$it.Top = 243.596142
$it.Width = 7.493701
$it.Height = 2.406220
$it = $g.GroupItems.Item(42)
$it.Left = 191.280236
$it.Top = 245.452441
$it.Width = 7.665591
$it.Height = 6.153071
$it = $g.GroupItems.Item(43)
$it.Left = 193.445906
$it.Top = 247.618031
$it.Width = 3.059370
$it.Height = 1.821811
$it = $g.GroupItems.Item(44)
$it.Left = 199.908425
$it.Top = 243.596142
$it.Width = 8.078110
$it.Height = 7.631181
$it = $g.GroupItems.Item(45)
$it.Left = 201.695906
$it.Top = 244.558661
$it.Width = 4.468740
$it.Height = 1.959370
$it = $g.GroupItems.Item(46)
$it.Left = 209.911496
$it.Top = 243.664882
$it.Width = 1.581181
$it.Height = 7.906220
$it = $g.GroupItems.Item(47)
$it.Left = 208.708425
$it.Top = 245.383701
$it.Width = 0.893701
$it.Height = 2.578110
$it = $g.GroupItems.Item(48)
$it.Left = 211.217795
$it.Top = 243.699291
$it.Width = 5.500000
$it.Height = 7.562441
$it = $g.GroupItems.Item(49)
$it.Left = 217.611496
$it.Top = 243.596142
$it.Width = 7.837480
$it.Height = 1.237480
$it = $g.GroupItems.Item(50)
$it.Left = 218.814646
$it.Top = 245.349291
$it.Width = 5.362441
$it.Height = 1.787480
$it = $g.GroupItems.Item(51)
$it.Left = 219.502126
$it.Top = 245.933701
$it.Width = 3.987480
$it.Height = 0.618740
$it = $g.GroupItems.Item(52)
$it.Left = 217.886535
$it.Top = 247.618031
$it.Width = 7.253071
$it.Height = 3.987480
$it = $g.GroupItems.Item(53)
$it.Left = 219.570866
$it.Top = 248.718031
$it.Width = 3.849921
$it.Height = 1.856220
$it = $g.GroupItems.Item(54)
$it.Left = 220.224016
$it.Top = 249.268031
$it.Width = 2.543701
$it.Height = 0.721811
$it = $g.GroupItems.Item(55)
$it.Left = 229.608425
$it.Top = 243.699291
$it.Width = 4.434331
$it.Height = 7.871811
$it = $g.GroupItems.Item(56)
$it.Left = 230.742756
$it.Top = 247.102441
$it.Width = 1.512441
$it.Height = 1.993701
$it = $g.GroupItems.Item(57)
$it.Left = 226.652126
$it.Top = 243.664882
$it.Width = 2.818740
$it.Height = 7.596850
$it = $g.GroupItems.Item(58)
$it.Left = 227.270866
$it.Top = 245.452441
$it.Width = 1.581181
$it.Height = 1.993701
$it = $g.GroupItems.Item(59)
$it.Left = 227.270866
$it.Top = 248.099291
$it.Width = 1.581181
$it.Height = 1.993701
$it = $g.GroupItems.Item(60)
$it.Left = 234.936535
$it.Top = 243.664882
$it.Width = 3.609370
$it.Height = 7.906220
$it = $g.GroupItems.Item(61)
$it.Left = 238.752126
$it.Top = 244.077402
$it.Width = 3.987480
$it.Height = 7.493701
$it = $g.GroupItems.Item(62)
$it.Left = 239.405276
$it.Top = 244.696142
$it.Width = 2.681181
$it.Height = 1.409370
$it = $g.GroupItems.Item(63)
$it.Left = 239.405276
$it.Top = 246.724252
$it.Width = 2.681181
$it.Height = 1.409370
$it = $g.GroupItems.Item(64)
$it.Left = 239.405276
$it.Top = 248.752441
$it.Width = 2.681181
$it.Height = 1.443701
$it = $g.GroupItems.Item(65)
$it.Left = 244.045906
$it.Top = 243.630551
$it.Width = 7.871811
$it.Height = 7.974961
$it = $g.GroupItems.Item(66)
$it.Left = 245.695906
$it.Top = 243.699291
$it.Width = 1.306220
$it.Height = 1.546850
$it = $g.GroupItems.Item(67)
$it.Left = 253.395906
$it.Top = 243.664882
$it.Width = 1.409370
$it.Height = 1.753071
$it = $g.GroupItems.Item(68)
$it.Left = 252.777165
$it.Top = 246.208661
$it.Width = 2.612441
$it.Height = 4.915591
$it = $g.GroupItems.Item(69)
$it.Left = 255.114646
$it.Top = 244.111811
$it.Width = 5.603071
$it.Height = 7.149921
$it = $g.GroupItems.Item(70)
$it.Left = 261.542756
$it.Top = 243.664882
$it.Width = 2.750000
$it.Height = 7.906220
$it = $g.GroupItems.Item(71)
$it.Left = 263.914646
$it.Top = 244.008661
$it.Width = 5.500000
$it.Height = 7.562441
$it = $g.GroupItems.Item(72)
$it.Left = 265.255276
$it.Top = 244.627402
$it.Width = 3.196850
$it.Height = 0.859370
$it = $g.GroupItems.Item(73)
$it.Left = 265.255276
$it.Top = 246.105512
$it.Width = 1.684331
$it.Height = 0.962441
$it = $g.GroupItems.Item(74)
$it.Left = 266.149055
$it.Top = 249.302441
$it.Width = 2.268740
$it.Height = 1.134331
$it = $g.GroupItems.Item(75)
$it.Left = 172.515827
$it.Top = 210.688583
$it.Width = 4.515984
$it.Height = 6.458189
$it = $g.GroupItems.Item(76)
$it.Left = 177.917008
$it.Top = 215.986614
$it.Width = 1.091339
$it.Height = 1.198819
$it = $g.GroupItems.Item(77)
$it.Left = 182.549055
$it.Top = 209.876457
$it.Width = 7.596850
$it.Height = 8.078110
$it = $g.GroupItems.Item(78)
$it.Left = 191.383386
$it.Top = 209.876457
$it.Width = 7.493701
$it.Height = 2.406220
$it = $g.GroupItems.Item(79)
$it.Left = 192.380236
$it.Top = 211.767087
$it.Width = 5.465591
$it.Height = 0.584331
$it = $g.GroupItems.Item(80)
$it.Left = 192.552126
$it.Top = 212.763937
$it.Width = 5.121811
$it.Height = 1.546850
$it = $g.GroupItems.Item(81)
$it.Left = 193.274016
$it.Top = 213.279606
$it.Width = 3.712441
$it.Height = 0.515591
$it = $g.GroupItems.Item(82)
$it.Left = 191.830236
$it.Top = 214.757717
$it.Width = 6.599921
$it.Height = 3.093701
$it = $g.GroupItems.Item(83)
$it.Left = 192.517795
$it.Top = 215.342047
$it.Width = 2.303071
$it.Height = 0.515591
$it = $g.GroupItems.Item(84)
$it.Left = 195.439685
$it.Top = 215.342047
$it.Width = 2.303071
$it.Height = 0.515591
$it = $g.GroupItems.Item(85)
$it.Left = 192.517795
$it.Top = 216.442047
$it.Width = 2.303071
$it.Height = 0.515591
$it = $g.GroupItems.Item(86)
$it.Left = 195.439685
$it.Top = 216.442047
$it.Width = 2.303071
$it.Height = 0.515591
$it = $g.GroupItems.Item(87)
$it.Left = 203.208425
$it.Top = 209.979606
$it.Width = 4.434331
$it.Height = 7.871811
$it = $g.GroupItems.Item(88)
$it.Left = 204.342756
$it.Top = 213.382677
$it.Width = 1.512441
$it.Height = 1.993701
$it = $g.GroupItems.Item(89)
$it.Left = 200.252126
$it.Top = 209.945197
$it.Width = 2.818740
$it.Height = 7.596850
$it = $g.GroupItems.Item(90)
$it.Left = 200.870866
$it.Top = 211.732677
$it.Width = 1.581181
$it.Height = 1.993701
$it = $g.GroupItems.Item(91)
$it.Left = 200.870866
$it.Top = 214.379606
$it.Width = 1.581181
$it.Height = 1.993701
$it = $g.GroupItems.Item(92)
$it.Left = 208.845906
$it.Top = 209.945197
$it.Width = 7.768740
$it.Height = 1.787480
$it = $g.GroupItems.Item(93)
$it.Left = 209.361496
$it.Top = 211.973307
$it.Width = 6.737480
$it.Height = 5.878110
$it = $g.GroupItems.Item(94)
$it.Left = 210.049055
$it.Top = 212.592047
$it.Width = 5.362441
$it.Height = 4.296850
$it = $g.GroupItems.Item(95)
$it.Left = 210.392756
$it.Top = 212.867087
$it.Width = 4.640551
$it.Height = 3.815591
$it = $g.GroupItems.Item(96)
$it.Left = 213.417795
$it.Top = 215.067087
$it.Width = 1.443701
$it.Height = 1.409370
$it = $g.GroupItems.Item(97)
$it.Left = 217.474016
$it.Top = 209.945197
$it.Width = 2.818740
$it.Height = 7.906220
$it = $g.GroupItems.Item(98)
$it.Left = 219.845906
$it.Top = 209.945197
$it.Width = 5.740551
$it.Height = 7.906220
$it = $g.GroupItems.Item(99)
$it.Left = 226.274016
$it.Top = 209.945197
$it.Width = 8.112441
$it.Height = 6.531181
$it = $g.GroupItems.Item(100)
$it.Left = 228.680236
$it.Top = 211.148346
$it.Width = 3.265591
$it.Height = 0.515591
$it = $g.GroupItems.Item(101)
$it.Left = 228.680236
$it.Top = 212.248346
$it.Width = 3.265591
$it.Height = 0.584331
$it = $g.GroupItems.Item(102)
$it.Left = 228.680236
$it.Top = 213.417087
$it.Width = 3.265591
$it.Height = 0.618740
$it = $g.GroupItems.Item(103)
$it.Left = 227.064646
$it.Top = 214.895197
$it.Width = 6.599921
$it.Height = 2.715591
$it = $g.GroupItems.Item(104)
$it.Left = 235.074016
$it.Top = 210.495197
$it.Width = 3.368740
$it.Height = 7.081181
$it = $g.GroupItems.Item(105)
$it.Left = 236.655276
$it.Top = 213.829606
$it.Width = 0.790551
$it.Height = 2.337480
$it = $g.GroupItems.Item(106)
$it.Left = 238.786535
$it.Top = 209.945197
$it.Width = 4.021811
$it.Height = 7.906220
$it = $g.GroupItems.Item(107)
$it.Left = 245.111496
$it.Top = 209.945197
$it.Width = 1.581181
$it.Height = 7.906220
$it = $g.GroupItems.Item(108)
$it.Left = 243.908425
$it.Top = 211.663937
$it.Width = 0.893701
$it.Height = 2.578110
$it = $g.GroupItems.Item(109)
$it.Left = 246.417795
$it.Top = 209.979606
$it.Width = 5.500000
$it.Height = 7.562441
$it = $g.GroupItems.Item(110)
$it.Left = 252.674016
$it.Top = 210.495197
$it.Width = 3.368740
$it.Height = 7.081181
$it = $g.GroupItems.Item(111)
$it.Left = 254.255276
$it.Top = 213.829606
$it.Width = 0.790551
$it.Height = 2.337480
$it = $g.GroupItems.Item(112)
$it.Left = 255.699055
$it.Top = 210.288976
$it.Width = 5.018740
$it.Height = 7.562441
$it = $g.GroupItems.Item(113)
$it.Left = 257.795906
$it.Top = 210.942047
$it.Width = 1.134331
$it.Height = 2.612441
$it = $g.GroupItems.Item(114)
$it.Left = 261.783386
$it.Top = 209.842047
$it.Width = 7.528110
$it.Height = 2.337480
$it = $g.GroupItems.Item(115)
$it.Left = 266.217795
$it.Top = 211.526457
$it.Width = 2.853071
$it.Height = 1.924961
$it = $g.GroupItems.Item(116)
$it.Left = 262.058425
$it.Top = 211.526457
$it.Width = 3.024961
$it.Height = 1.753071
$it = $g.GroupItems.Item(117)
$it.Left = 261.680236
$it.Top = 212.798346
$it.Width = 7.906220
$it.Height = 5.087480
$it = $g.GroupItems.Item(118)
$it.Left = 172.674803
$it.Top = 172.307087
$it.Width = 4.124961
$it.Height = 6.281969
$it = $g.GroupItems.Item(119)
$it.Left = 173.474016
$it.Top = 172.990315
$it.Width = 2.333150
$it.Height = 1.980787
$it = $g.GroupItems.Item(120)
$it.Left = 173.474016
$it.Top = 175.641417
$it.Width = 2.483543
$it.Height = 2.268740
$it = $g.GroupItems.Item(121)
$it.Left = 177.917008
$it.Top = 177.514882
$it.Width = 1.091339
$it.Height = 1.198819
$it = $g.GroupItems.Item(122)
$it.Left = 182.239685
$it.Top = 171.473465
$it.Width = 3.024961
$it.Height = 7.906220
$it = $g.GroupItems.Item(123)
$it.Left = 185.470866
$it.Top = 171.439134
$it.Width = 4.778110
$it.Height = 3.334331
$it = $g.GroupItems.Item(124)
$it.Left = 185.505276
$it.Top = 175.014094
$it.Width = 4.778110
$it.Height = 4.365591
$it = $g.GroupItems.Item(125)
$it.Left = 186.158425
$it.Top = 176.767244
$it.Width = 1.237480
$it.Height = 1.512441
$it = $g.GroupItems.Item(126)
$it.Left = 191.452126
$it.Top = 171.404724
$it.Width = 7.356220
$it.Height = 2.268740
$it = $g.GroupItems.Item(127)
$it.Left = 191.211496
$it.Top = 173.914094
$it.Width = 7.906220
$it.Height = 5.431181
$it = $g.GroupItems.Item(128)
$it.Left = 199.839685
$it.Top = 171.404724
$it.Width = 7.940551
$it.Height = 8.009370
$it = $g.GroupItems.Item(129)
$it.Left = 199.977165
$it.Top = 172.951575
$it.Width = 1.134331
$it.Height = 1.787480
$it = $g.GroupItems.Item(130)
$it.Left = 202.005276
$it.Top = 173.123465
$it.Width = 5.912441
$it.Height = 6.290551
$it = $g.GroupItems.Item(131)
$it.Left = 208.605276
$it.Top = 171.370394
$it.Width = 8.078110
$it.Height = 8.043701
$it = $g.GroupItems.Item(132)
$it.Left = 208.708425
$it.Top = 172.917244
$it.Width = 1.100000
$it.Height = 1.718740
$it = $g.GroupItems.Item(133)
$it.Left = 210.977165
$it.Top = 173.535984
$it.Width = 5.568740
$it.Height = 5.740551
$it = $g.GroupItems.Item(134)
$it.Left = 220.808425
$it.Top = 171.507874
$it.Width = 4.434331
$it.Height = 7.871811
$it = $g.GroupItems.Item(135)
$it.Left = 221.942756
$it.Top = 174.910945
$it.Width = 1.512441
$it.Height = 1.993701
$it = $g.GroupItems.Item(136)
$it.Left = 217.852126
$it.Top = 171.473465
$it.Width = 2.818740
$it.Height = 7.596850
$it = $g.GroupItems.Item(137)
$it.Left = 218.470866
$it.Top = 173.260945
$it.Width = 1.581181
$it.Height = 1.993701
$it = $g.GroupItems.Item(138)
$it.Left = 218.470866
$it.Top = 175.907874
$it.Width = 1.581181
$it.Height = 1.993701
$it = $g.GroupItems.Item(139)
$it.Left = 226.411496
$it.Top = 171.404724
$it.Width = 7.837480
$it.Height = 1.237480
$it = $g.GroupItems.Item(140)
$it.Left = 227.614646
$it.Top = 173.157874
$it.Width = 5.362441
$it.Height = 1.787480
$it = $g.GroupItems.Item(141)
$it.Left = 228.302126
$it.Top = 173.742205
$it.Width = 3.987480
$it.Height = 0.618740
$it = $g.GroupItems.Item(142)
$it.Left = 226.686535
$it.Top = 175.426614
$it.Width = 7.253071
$it.Height = 3.987480
$it = $g.GroupItems.Item(143)
$it.Left = 228.370866
$it.Top = 176.526614
$it.Width = 3.849921
$it.Height = 1.856220
$it = $g.GroupItems.Item(144)
$it.Left = 229.024016
$it.Top = 177.076614
$it.Width = 2.543701
$it.Height = 0.721811
$it = $g.GroupItems.Item(145)
$it.Left = 235.177165
$it.Top = 171.542205
$it.Width = 7.631181
$it.Height = 7.871811
$it = $g.GroupItems.Item(146)
$it.Left = 236.586535
$it.Top = 172.607874
$it.Width = 6.428110
$it.Height = 5.121811
$it = $g.GroupItems.Item(147)
$it.Left = 236.724016
$it.Top = 175.942205
$it.Width = 3.506220
$it.Height = 3.506220
$it = $g.GroupItems.Item(148)
$it.Left = 240.127165
$it.Top = 177.592205
$it.Width = 2.750000
$it.Height = 1.856220
$it = $g.GroupItems.Item(149)
$it.Left = 245.008425
$it.Top = 171.645354
$it.Width = 5.774961
$it.Height = 2.303071
$it = $g.GroupItems.Item(150)
$it.Left = 245.695906
$it.Top = 172.160945
$it.Width = 4.399921
$it.Height = 0.378110
$it = $g.GroupItems.Item(151)
$it.Left = 245.695906
$it.Top = 173.020394
$it.Width = 4.399921
$it.Height = 0.378110
$it = $g.GroupItems.Item(152)
$it.Left = 244.114646
$it.Top = 174.292205
$it.Width = 7.596850
$it.Height = 0.515591
$it = $g.GroupItems.Item(153)
$it.Left = 243.942756
$it.Top = 175.151575
$it.Width = 7.940551
$it.Height = 4.056220
$it = $g.GroupItems.Item(154)
$it.Left = 245.661496
$it.Top = 175.701575
$it.Width = 1.890551
$it.Height = 0.378110
$it = $g.GroupItems.Item(155)
$it.Left = 248.274016
$it.Top = 175.701575
$it.Width = 1.890551
$it.Height = 0.378110
$it = $g.GroupItems.Item(156)
$it.Left = 245.661496
$it.Top = 176.526614
$it.Width = 1.890551
$it.Height = 0.343701
$it = $g.GroupItems.Item(157)
$it.Left = 248.274016
$it.Top = 176.526614
$it.Width = 1.890551
$it.Height = 0.343701
$it = $g.GroupItems.Item(158)
$it.Left = 254.599055
$it.Top = 171.610945
$it.Width = 0.653071
$it.Height = 7.631181
$it = $g.GroupItems.Item(159)
$it.Left = 253.086535
$it.Top = 172.332835
$it.Width = 0.687480
$it.Height = 5.671811
$it = $g.GroupItems.Item(160)
$it.Left = 255.492756
$it.Top = 171.473465
$it.Width = 5.018740
$it.Height = 3.609370
$it = $g.GroupItems.Item(161)
$it.Left = 257.727165
$it.Top = 173.707874
$it.Width = 1.443701
$it.Height = 1.478110
$it = $g.GroupItems.Item(162)
$it.Left = 256.111496
$it.Top = 175.357874
$it.Width = 4.193701
$it.Height = 4.056220
$it = $g.GroupItems.Item(163)
$it.Left = 256.799055
$it.Top = 176.010945
$it.Width = 1.031181
$it.Height = 1.924961
$it = $g.GroupItems.Item(164)
$it.Left = 258.517795
$it.Top = 176.010945
$it.Width = 1.100000
$it.Height = 1.924961
$it = $g.GroupItems.Item(165)
$it.Left = 261.474016
$it.Top = 171.439134
$it.Width = 7.803071
$it.Height = 7.940551
$it = $g.GroupItems.Item(166)
$it.Left = 263.055276
$it.Top = 173.260945
$it.Width = 6.256220
$it.Height = 6.049921
$it = $g.GroupItems.Item(167)
$it.Left = 252.674016
$it.Top = 181.527480
$it.Width = 3.368740
$it.Height = 7.081181
$it = $g.GroupItems.Item(168)
$it.Left = 254.255276
$it.Top = 184.861890
$it.Width = 0.790551
$it.Height = 2.337480
$it = $g.GroupItems.Item(169)
$it.Left = 255.699055
$it.Top = 181.321260
$it.Width = 5.018740
$it.Height = 7.562441
$it = $g.GroupItems.Item(170)
$it.Left = 257.795906
$it.Top = 181.974331
$it.Width = 1.134331
$it.Height = 2.612441
$it = $g.GroupItems.Item(171)
$it.Left = 261.783386
$it.Top = 180.874331
$it.Width = 7.528110
$it.Height = 2.337480
$it = $g.GroupItems.Item(172)
$it.Left = 266.217795
$it.Top = 182.558740
$it.Width = 2.853071
$it.Height = 1.924961
$it = $g.GroupItems.Item(173)
$it.Left = 262.058425
$it.Top = 182.558740
$it.Width = 3.024961
$it.Height = 1.753071
$it = $g.GroupItems.Item(174)
$it.Left = 261.680236
$it.Top = 183.830630
$it.Width = 7.906220
$it.Height = 5.087480
$it = $g.GroupItems.Item(175)
$it.Left = 180.735748
$it.Top = 143.313543
$it.Width = 5.349606
$it.Height = 6.307795
$it = $g.GroupItems.Item(176)
$it.Left = 182.497480
$it.Top = 144.172992
$it.Width = 1.830394
$it.Height = 2.792913
$it = $g.GroupItems.Item(177)
$it.Left = 186.717008
$it.Top = 148.547165
$it.Width = 1.091339
$it.Height = 1.198819
$it = $g.GroupItems.Item(178)
$it.Left = 191.177165
$it.Top = 142.402677
$it.Width = 7.974961
$it.Height = 2.818740
$it = $g.GroupItems.Item(179)
$it.Left = 192.449055
$it.Top = 144.052677
$it.Width = 5.396850
$it.Height = 6.428110
$it = $g.GroupItems.Item(180)
$it.Left = 193.136535
$it.Top = 145.530787
$it.Width = 4.021811
$it.Height = 0.618740
$it = $g.GroupItems.Item(181)
$it.Left = 193.136535
$it.Top = 146.768268
$it.Width = 4.021811
$it.Height = 0.618740
$it = $g.GroupItems.Item(182)
$it.Left = 194.167795
$it.Top = 148.005748
$it.Width = 4.846850
$it.Height = 2.474961
$it = $g.GroupItems.Item(183)
$it.Left = 201.008425
$it.Top = 142.849528
$it.Width = 5.843701
$it.Height = 2.990551
$it = $g.GroupItems.Item(184)
$it.Left = 201.695906
$it.Top = 143.502677
$it.Width = 4.503071
$it.Height = 1.684331
$it = $g.GroupItems.Item(185)
$it.Left = 200.286535
$it.Top = 146.527638
$it.Width = 3.265591
$it.Height = 3.884331
$it = $g.GroupItems.Item(186)
$it.Left = 200.939685
$it.Top = 147.146378
$it.Width = 1.959370
$it.Height = 2.131181
$it = $g.GroupItems.Item(187)
$it.Left = 204.308425
$it.Top = 146.527638
$it.Width = 3.231181
$it.Height = 3.884331
$it = $g.GroupItems.Item(188)
$it.Left = 204.961496
$it.Top = 147.146378
$it.Width = 1.924961
$it.Height = 2.131181
$it = $g.GroupItems.Item(189)
$it.Left = 208.742756
$it.Top = 142.505748
$it.Width = 2.750000
$it.Height = 7.906220
$it = $g.GroupItems.Item(190)
$it.Left = 211.252126
$it.Top = 142.471417
$it.Width = 5.465591
$it.Height = 7.940551
$it = $g.GroupItems.Item(191)
$it.Left = 212.764646
$it.Top = 144.568268
$it.Width = 1.478110
$it.Height = 1.134331
$it = $g.GroupItems.Item(192)
$it.Left = 212.764646
$it.Top = 146.287008
$it.Width = 1.478110
$it.Height = 1.168740
$it = $g.GroupItems.Item(193)
$it.Left = 212.764646
$it.Top = 148.074488
$it.Width = 1.478110
$it.Height = 1.134331
$it = $g.GroupItems.Item(194)
$it.Left = 217.439685
$it.Top = 142.471417
$it.Width = 7.837480
$it.Height = 7.974961
$it = $g.GroupItems.Item(195)
$it.Left = 226.308425
$it.Top = 142.918268
$it.Width = 7.734331
$it.Height = 7.493701
$it = $g.GroupItems.Item(196)
$it.Left = 230.983386
$it.Top = 145.083858
$it.Width = 3.368740
$it.Height = 3.334331
$it = $g.GroupItems.Item(197)
$it.Left = 235.314646
$it.Top = 142.643228
$it.Width = 7.631181
$it.Height = 3.643701
$it = $g.GroupItems.Item(198)
$it.Left = 236.345906
$it.Top = 144.740157
$it.Width = 2.096850
$it.Height = 0.549921
$it = $g.GroupItems.Item(199)
$it.Left = 239.817795
$it.Top = 144.740157
$it.Width = 2.096850
$it.Height = 0.549921
$it = $g.GroupItems.Item(200)
$it.Left = 236.311496
$it.Top = 145.702677
$it.Width = 2.096850
$it.Height = 0.515591
$it = $g.GroupItems.Item(201)
$it.Left = 239.783386
$it.Top = 145.702677
$it.Width = 2.096850
$it.Height = 0.515591
$it = $g.GroupItems.Item(202)
$it.Left = 235.486535
$it.Top = 146.665118
$it.Width = 7.287480
$it.Height = 3.781181
$it = $g.GroupItems.Item(203)
$it.Left = 243.908425
$it.Top = 146.458898
$it.Width = 7.940551
$it.Height = 3.987480
$it = $g.GroupItems.Item(204)
$it.Left = 246.314646
$it.Top = 147.662047
$it.Width = 3.162441
$it.Height = 1.271811
$it = $g.GroupItems.Item(205)
$it.Left = 244.217795
$it.Top = 142.746378
$it.Width = 7.424961
$it.Height = 3.506220
$it = $g.GroupItems.Item(206)
$it.Left = 247.277165
$it.Top = 143.399528
$it.Width = 1.374961
$it.Height = 0.687480
$it = $g.GroupItems.Item(207)
$it.Left = 245.317795
$it.Top = 144.671417
$it.Width = 1.340551
$it.Height = 0.996850
$it = $g.GroupItems.Item(208)
$it.Left = 247.277165
$it.Top = 144.671417
$it.Width = 1.374961
$it.Height = 0.996850
$it = $g.GroupItems.Item(209)
$it.Left = 249.270866
$it.Top = 144.671417
$it.Width = 1.271811
$it.Height = 0.996850
$it = $g.GroupItems.Item(210)
$it.Left = 253.395906
$it.Top = 142.505748
$it.Width = 1.409370
$it.Height = 1.753071
$it = $g.GroupItems.Item(211)
$it.Left = 252.777165
$it.Top = 145.049528
$it.Width = 2.612441
$it.Height = 4.915591
$it = $g.GroupItems.Item(212)
$it.Left = 255.114646
$it.Top = 142.952677
$it.Width = 5.603071
$it.Height = 7.149921
$it = $g.GroupItems.Item(213)
$it.Left = 261.542756
$it.Top = 142.505748
$it.Width = 2.750000
$it.Height = 7.906220
$it = $g.GroupItems.Item(214)
$it.Left = 263.914646
$it.Top = 142.849528
$it.Width = 5.500000
$it.Height = 7.562441
$it = $g.GroupItems.Item(215)
$it.Left = 265.255276
$it.Top = 143.468268
$it.Width = 3.196850
$it.Height = 0.859370
$it = $g.GroupItems.Item(216)
$it.Left = 265.255276
$it.Top = 144.946378
$it.Width = 1.684331
$it.Height = 0.962441
$it = $g.GroupItems.Item(217)
$it.Left = 266.149055
$it.Top = 148.143228
$it.Width = 2.268740
$it.Height = 1.134331
$it = $g.GroupItems.Item(218)
$it.Left = 269.930236
$it.Top = 247.753386
$it.Width = 2.739685
$it.Height = 0.000000
$it = $g.GroupItems.Item(219)
$it.Left = 269.930236
$it.Top = 214.033701
$it.Width = 2.739685
$it.Height = 0.000000
$it = $g.GroupItems.Item(220)
$it.Left = 269.930236
$it.Top = 180.313937
$it.Width = 2.739685
$it.Height = 0.000000
$it = $g.GroupItems.Item(221)
$it.Left = 269.930236
$it.Top = 146.594252
$it.Width = 2.739685
$it.Height = 0.000000
$it = $g.GroupItems.Item(222)
$it.Left = 272.670000
$it.Top = 267.985197
$it.Width = 0.000000
$it.Height = 2.739685
$it = $g.GroupItems.Item(223)
$it.Left = 314.613701
$it.Top = 267.985197
$it.Width = 0.000000
$it.Height = 2.739685
$it = $g.GroupItems.Item(224)
$it.Left = 356.557402
$it.Top = 267.985197
$it.Width = 0.000000
$it.Height = 2.739685
$it = $g.GroupItems.Item(225)
$it.Left = 398.501024
$it.Top = 267.985197
$it.Width = 0.000000
$it.Height = 2.739685
$it = $g.GroupItems.Item(226)
$it.Left = 270.566693
$it.Top = 272.826457
$it.Width = 4.206614
$it.Height = 6.230394
$it = $g.GroupItems.Item(227)
$it.Left = 271.348740
$it.Top = 273.453858
$it.Width = 2.638268
$it.Height = 4.971417
$it = $g.GroupItems.Item(228)
$it.Left = 307.715039
$it.Top = 272.826457
$it.Width = 4.008976
$it.Height = 6.144488
$it = $g.GroupItems.Item(229)
$it.Left = 312.518976
$it.Top = 272.916693
$it.Width = 4.172205
$it.Height = 6.140157
$it = $g.GroupItems.Item(230)
$it.Left = 317.404488
$it.Top = 272.826457
$it.Width = 4.206614
$it.Height = 6.230394
$it = $g.GroupItems.Item(231)
$it.Left = 318.186535
$it.Top = 273.453858
$it.Width = 2.638268
$it.Height = 4.971417
$it = $g.GroupItems.Item(232)
$it.Left = 349.568504
$it.Top = 272.916693
$it.Width = 4.172205
$it.Height = 6.140157
$it = $g.GroupItems.Item(233)
$it.Left = 354.454016
$it.Top = 272.826457
$it.Width = 4.206614
$it.Height = 6.230394
$it = $g.GroupItems.Item(234)
$it.Left = 355.236063
$it.Top = 273.453858
$it.Width = 2.638268
$it.Height = 4.971417
$it = $g.GroupItems.Item(235)
$it.Left = 359.348189
$it.Top = 272.826457
$it.Width = 4.206614
$it.Height = 6.230394
$it = $g.GroupItems.Item(236)
$it.Left = 360.130236
$it.Top = 273.453858
$it.Width = 2.638268
$it.Height = 4.971417
$it = $g.GroupItems.Item(237)
$it.Left = 391.611024
$it.Top = 272.916693
$it.Width = 4.000315
$it.Height = 6.054252
$it = $g.GroupItems.Item(238)
$it.Left = 396.406299
$it.Top = 272.916693
$it.Width = 4.172205
$it.Height = 6.140157
$it = $g.GroupItems.Item(239)
$it.Left = 401.291890
$it.Top = 272.826457
$it.Width = 4.206614
$it.Height = 6.230394
$it = $g.GroupItems.Item(240)
$it.Left = 402.073937
$it.Top = 273.453858
$it.Width = 2.638268
$it.Height = 4.971417